$d = $word.ActiveDocument

$old = "Implemented parallel & asynchronous programming which vastly improved the performance during word combine & pdf conversion."
$new = "Implemented parallel & asynchronous programming to download/upload documents from/to Sharepoint & S3 which vastly improved the performance during word to pdf conversion."

$found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Target text for replacement was not found."
}

Write-Output "Replaced paragraph text. Found=$found"
